# pet-simulator Upgrade_升级属性表: add pet bagMaxCapacity 100
# - Column F (PetNum / 宠物背包增加数量) changes from a single "int" to a
#   list "int[]" of per-level values.
# - Row 9 (Page_Title_8) Diamond-cost list (D9) is replaced with the
#   standard 5-value progression, and the PetNum cell (F9) becomes the
#   per-level capacity list " 5|5|10|20|30" (note leading space, matches
#   source data) instead of the flat literal 5.
# - F1/F2/F9 pick up the alternate header/data font style already used
#   by D9 (style index 4) via PasteSpecial formats from D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content updates -------------------------------------------------
$ws.Range("F1").Value = "int[]"
$ws.Range("D9").Value = "10000|40000|200000|1000000|5000000"
$ws.Range("F9").Value = " 5|5|10|20|30"

# --- style updates -----------------------------------------------------
# F1, F2 and F9 switch from style 1 to the style already used by D9
# (style 4). Copy D9's format (which is unchanged) onto them.
$ws.Range("D9").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F2").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F9").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# --- selection / view ---------------------------------------------------
$ws.Range("E9").Select()
